$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stuff-descriptor")
$ws.Cells.Item(1, 16).Value = "value"
